$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right after
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. At the end of the document, split the final (italic, image-prompt)
#    paragraph so that a new paragraph is inserted right before it. That new
#    paragraph carries the bold "Play Fruit Million Free..." heading text
#    that used to live in the meta-description paragraph we just removed.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastRange = $d.Range($lastStart, $lastEnd)
$lastRange.InsertParagraphBefore()

# The freshly-inserted (currently empty) paragraph is now second-to-last.
$newHeadingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$headingStart = $newHeadingPara.Range.Start
$headingEnd = $newHeadingPara.Range.End
$headingRange = $d.Range($headingStart, $headingEnd)
$headingRange.Text = "Play Fruit Million Free: Impressive Graphics and Expanding Wilds"
$headingRange.Font.Bold = 1
$headingRange.Font.Italic = 0

# ---------------------------------------------------------------------------
# 3. Replace the text of the last paragraph (still italic) with the
#    meta-description body copy.
# ---------------------------------------------------------------------------
$imagePromptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptStart = $imagePromptPara.Range.Start
$promptEnd = $imagePromptPara.Range.End
$promptRange = $d.Range($promptStart, $promptEnd)
$promptRange.Text = "Read our review of Fruit Million by BGaming. Play for free and experience the game's unique visuals, expanding Wilds, and impressive RTP."
